$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates from the coin price/volume refresh (19 Jan 2023 14:33 UTC run).
# Values are stored as text (not numbers) to match the source data feed formatting,
# so each cell is forced to Text format before the write and restored to the default
# style afterward (only the value should differ from the original workbook).
$cellUpdates = @{
    'D2' = '288.76'
    'E2' = '-4.52%'
    'D3' = '30.19'
    'E3' = '-7.16%'
    'D4' = '4.894'
    'E4' = '-3.02%'
    'D5' = '0.07237'
    'E5' = '-6.18%'
    'D6' = '1.769'
    'E6' = '-14.98%'
    'D7' = '7.572'
    'E8' = '-2.13%'
    'D9' = '0.9012'
    'E9' = '-2.65%'
    'D10' = '0.1665'
    'E10' = '-5.31%'
    'D11' = '0.07821'
    'E11' = '-1.82%'
    'D12' = '0.07949'
    'E12' = '-6.32%'
    'D13' = '0.03048'
    'E13' = '-0.07%'
    'E14' = '0.39%'
    'E15' = '-1.90%'
    'D16' = '0.005781'
    'E16' = '0.26%'
    'D17' = '3.470'
    'E17' = '0.13%'
    'E18' = '-3.80%'
    'D19' = '0.3318'
    'E19' = '-1.19%'
    'D20' = '0.1294'
    'E20' = '-0.43%'
    'D21' = '3.966'
    'E21' = '-9.45%'
    'D22' = '0.2179'
    'E22' = '10.36%'
    'D23' = '0.04492'
    'E23' = '-0.76%'
    'D24' = '0.001212'
    'E24' = '-2.25%'
    'D25' = '0.004419'
    'E25' = '6.49%'
    'D26' = '0.0001300'
    'E26' = '4.05%'
    'D27' = '0.0003394'
    'D39' = '0.01557'
    'E39' = '-9.60%'
    'D40' = '0.04322'
    'E40' = '-8.14%'
    'D41' = '0.007365'
    'E41' = '-1.43%'
    'D42' = '0.01008'
    'D43' = '0.1305'
    'E43' = '-4.19%'
    'D44' = '0.002005'
    'E44' = '-10.25%'
    'D45' = '0.009337'
    'E45' = '-11.86%'
    'D46' = '0.00005881'
    'E46' = '-5.25%'
    'D47' = '0.00000000750'
    'E47' = '0.05%'
    'E48' = '259.07%'
    'D49' = '0.002898'
    'E49' = '-3.34%'
    'D50' = '0.00002100'
    'E50' = '0.05%'
    'D51' = '0.0002000'
    'E51' = '0.05%'
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.Style = "Normal"
}
